# Update cryptocurrency price/volume data per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.538.01'
$ws.Range("E2").Value = '  -4.45%  '
$ws.Range("D3").Value = '2.492.84'
$ws.Range("E3").Value = '  -5.53%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'542.80"
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").Value = "'146.85"
$ws.Range("E6").Value = '  -4.59%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("D9").Value = '2.520.63'
$ws.Range("E9").Value = '  -4.55%  '
$ws.Range("E10").Value = '  -2.86%  '
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").Value = "'5.49"
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").Value = '2.940.53'
$ws.Range("E14").Value = '  -5.34%  '
$ws.Range("E15").Value = '  -4.26%  '
$ws.Range("D16").Value = '59.595.80'
$ws.Range("E16").Value = '  -4.24%  '
$ws.Range("E17").Value = '  -2.29%  '
$ws.Range("D18").Value = '2.504.27'
$ws.Range("E18").Value = '  -5.15%  '
$ws.Range("D19").Value = "'11.41"
$ws.Range("E19").Value = '  -2.13%  '
$ws.Range("E20").Value = '  -3.25%  '
$ws.Range("D21").Value = "'326.92"
$ws.Range("E21").Value = '  -3.49%  '
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("E23").Value = '  -4.02%  '
$ws.Range("D24").Value = "'61.42"
$ws.Range("E24").Value = '  -1.88%  '
$ws.Range("D25").Value = "'0.445"
$ws.Range("E25").Value = '  -10.88%  '
$ws.Range("E26").Value = '  +0.89%  '
$ws.Range("E27").Value = '  -3.15%  '
$ws.Range("D28").Value = "'7.83"
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("D29").Value = '0.0₃0796'
$ws.Range("E29").Value = '  -4.42%  '
$ws.Range("D30").Value = "'1.29"
$ws.Range("E30").Value = '  -4.83%  '
$ws.Range("E31").Value = '  -3.84%  '
$ws.Range("E32").Value = '  -2.62%  '
$ws.Range("D33").Value = "'0.997"
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("D34").Value = "'158.27"
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").Value = '  +1.59%  '
$ws.Range("D36").Value = "'19.07"
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").Value = "'4.48"
$ws.Range("E37").Value = '  -5.37%  '
$ws.Range("E38").Value = '  +0.76%  '
$ws.Range("D39").Value = "'5.98"
$ws.Range("E39").Value = '  -1.79%  '
$ws.Range("D40").Value = "'314.48"
$ws.Range("E40").Value = '  -5.46%  '
$ws.Range("D41").Value = "'36.76"
$ws.Range("E41").Value = '  -3.17%  '
$ws.Range("E42").Value = '  -3.38%  '
$ws.Range("D43").Value = "'0.829"
$ws.Range("E43").Value = '  -7.89%  '
$ws.Range("D44").Value = "'0.994"
$ws.Range("E44").Value = '  -0.46%  '
$ws.Range("D45").Value = "'0.605"
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'127.16"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").Value = "'10.71"
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("E48").Value = '  -2.04%  '
$ws.Range("D49").Value = "'0.0941"
$ws.Range("E49").Value = '  -1.76%  '
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").Value = "'18.73"
$ws.Range("E51").Value = '  -4.86%  '
